$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handed back"
$overview.Range("C3").Value = "Handed back"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Handed back"
$zhcn.Range("G3").Value = "2016-01-08 15:26:54"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Handed back"
$dede.Range("G3").Value = "2016-01-08 15:27:17"
